$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '50.908.72'
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").Value = '2.938.28'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '378.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.96%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -1.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.99'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").Value = '3.395.78'
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.03%  '
$ws.Range("B15").Value = 'Uniswap'
$ws.Range("C15").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '12.18'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +69.76%  '
$ws.Range("E16").Value = '  +3.29%  '
$ws.Range("D17").Value = '2.944.88'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.994'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("D19").Value = '50.877.41'
$ws.Range("E19").Value = '  -0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.43%  '
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.95%  '
$ws.Range("E25").Value = '  +13.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.34%  '
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.58%  '
$ws.Range("E30").Value = '  -3.23%  '
$ws.Range("E31").Value = '  -2.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.41'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("E35").Value = '  -0.94%  '
$ws.Range("E36").Value = '  -2.73%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  +3.47%  '
$ws.Range("E39").Value = '  +1.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.53'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.27%  '
$ws.Range("E41").Value = '  +1.78%  '
$ws.Range("E42").Value = '  -3.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '119.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("E45").Value = '  +7.02%  '
$ws.Range("E46").Value = '  -1.88%  '
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("D48").Value = '2.008.99'
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("E49").Value = '  -4.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0314'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.11%  '
$ws.Range("E51").Value = '  +4.48%  '
